$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Row 12 (station 4 / CTD): rename "CTD-ST1" -> 1, and tag Type column with "CTD"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "CTD"

# Row 13 (station 5 / CTD): rename "CTD-ST2" -> 2, and tag Type column with "CTD"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "CTD"

# Row 14 (station 6 / CTD): rename "CTD-ST3" -> 3, and tag Type column with "CTD"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "CTD"

# Update the active selection to E2 (matches the saved view in the target file)
$ws.Range("E2").Select()
